$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq "Administrator, Miss Dina Nasr") {
        $cell.Value = "Miss Dina Nasr, Administrator"
    }
}
